# Flip the sign of the plotted loading values so the figure aesthetics
# (axis orientation / component sign) match the updated toy example output.

$wb = $excel.ActiveWorkbook

# --- Sheet "A.hat": columns B (Component 1) and C (Component 2), rows 2-61 ---
$wsA = $wb.Worksheets.Item("A.hat")
for ($r = 2; $r -le 61; $r++) {
    foreach ($col in @("B", "C")) {
        $cell = $wsA.Range("$col$r")
        $cell.Value = -1 * $cell.Value2
    }
}

# --- Sheet "Phi.hat": columns A (Component 1) and B (Component 2), rows 2-102 ---
$wsPhi = $wb.Worksheets.Item("Phi.hat")
for ($r = 2; $r -le 102; $r++) {
    foreach ($col in @("A", "B")) {
        $cell = $wsPhi.Range("$col$r")
        $cell.Value = -1 * $cell.Value2
    }
}
